$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range
$x = $r.XML()
Write-Output $x
